$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Database tables" section appended below the existing component table.
$ws.Range("A22").Value = "Database tables"

# Data rows for "admins" and "posts" tables
$ws.Range("A24").Value = "admins"
$ws.Range("A25").Value = "posts"
$ws.Range("B24").Value = "`$uid"
$ws.Range("B25").Value = "`$parent->`$child"
$ws.Range("C24").Value = "scope"
$ws.Range("C25").Value = "author, text, replyTo, threadId, timestamp"

# Column headers for the new table
$ws.Range("A23").Value = "table name"
$ws.Range("B23").Value = "key structure"
$ws.Range("C23").Value = "notable fields"

# Data row for "postvotes" table
$ws.Range("A26").Value = "postvotes"
$ws.Range("B26").Value = "`$parent->`$child"
$ws.Range("C26").Value = "up,down,timestamp"

# Data row for "profile" table (key column only for now)
$ws.Range("A27").Value = "profile"
$ws.Range("B27").Value = "`$uid"

# Data rows for "threads" and "uservotes" tables (name + key columns)
$ws.Range("A28").Value = "threads"
$ws.Range("B28").Value = "`$threadId"
$ws.Range("A29").Value = "uservotes"
$ws.Range("B29").Value = "`$uid->`$parent-`$child"

# Remaining notable-fields column
$ws.Range("C27").Value = "uid, username, email, picture_link"
$ws.Range("C28").Value = "owner, created, title, href, allowAnonymous"
$ws.Range("C29").Value = "up,down,timestamp"

$ws.Range("B25").Select()
